$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Div Original Name" column (column D), shifting subsequent
# columns left by one, matching the commit's base-update that dropped
# this header/field from the sheet.
$ws.Range("D:D").Delete()
